# daily auto push: 2026-01-13 06:50 UTC
# Insert one new data row just above the old row 627 (2026/12/29 ...),
# which pushes the existing rows 627-668 down to 628-669 and grows the
# used range from A1:D668 to A1:D669. Then populate the newly inserted
# row with the new observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 627..668 down to 628..669 by inserting a blank row at 627.
$ws.Rows.Item(627).Insert()

# Fill the new row 627. The date column holds plain text like "2026/12/29"
# elsewhere in the sheet (not a real date), so prefix with an apostrophe to
# stop Excel's automatic date parsing, then clear the resulting formatting
# so the cell ends up with no explicit style - matching its neighbours.
$ws.Range("A627").Value = "'2026/01/13"
$ws.Range("A627").ClearFormats()

$ws.Range("B627").Value = "火"
$ws.Range("C627").Value = 13
$ws.Range("D627").Value = 201
